# Updated CHE model - 2025-08-15 15:09
$wb = $excel.ActiveWorkbook

# --- 1. Rename commodity "elc_sol-CHE" -> "elc_spv-CHE" on the re_profiles sheet (column E, rows 4-123)
$wsProfiles = $wb.Worksheets.Item("re_profiles")
$rng = $wsProfiles.Range("E4:E123")
for ($i = 1; $i -le $rng.Rows.Count; $i++) {
    $cell = $rng.Cells.Item($i, 1)
    if ($cell.Value2 -eq "elc_sol-CHE") {
        $cell.Value2 = "elc_spv-CHE"
    }
}

# --- 2. Re-order the two comma-separated timeslice lists on ev_charging_uc (C13/C14)
#        G7/G8 on the same sheet already reference them via formulas (=C14 / =C13)
#        so they will pick up the new order automatically on recalculation.
$wsUc = $wb.Worksheets.Item("ev_charging_uc")

$wsUc.Range("C13").Value2 = "S1b0205h08,S1b0205h10,S2c0415h14,S4aH2,S1b0205h07,S2c0415h10,S2c0415h16,S2c0415h18,S2d0427h12,S3aH3,S1aH6,S2c0415h12,S1b0205h16,S2aH6,S2d0427h09,S3aH5,S5aH3,S1aH5,S1b0205h09,S1b0205h14,S1b0205h15,S2aH4,S2aH5,S2c0415h08,S5aH2,S5aH5,S6aH2,S1b0205h18,S2aH7,S2d0427h08,S3aH4,S4aH3,S4aH4,S2d0427h07,S2d0427h10,S2d0427h17,S6aH5,S1aH4,S2aH2,S2d0427h15,S4aH6,S2c0415h07,S2c0415h13,S2c0415h15,S2c0415h17,S2d0427h11,S6aH7,S1aH3,S1b0205h13,S2d0427h13,S4aH5,S1b0205h17,S2d0427h14,S2d0427h16,S5aH4,S6aH3,S6aH4,S1b0205h12,S2c0415h09,S2c0415h11,S3aH7,S1aH2,S1aH7,S1b0205h11,S2aH3,S2d0427h18,S3aH6,S4aH7,S5aH6,S5aH7,S6aH6,S3aH2"

$wsUc.Range("C14").Value2 = "S1b0205h19,S2aH1,S2c0415h01,S5aH8,S1b0205h23,S2c0415h06,S2d0427h04,S2d0427h24,S2c0415h23,S6aH8,S1aH8,S1b0205h22,S2d0427h02,S2d0427h06,S2d0427h23,S5aH1,S2c0415h20,S1b0205h02,S2d0427h01,S2d0427h19,S2c0415h21,S4aH8,S2c0415h04,S2d0427h21,S4aH1,S6aH1,S1aH1,S1b0205h04,S1b0205h06,S1b0205h21,S2c0415h02,S2c0415h19,S2d0427h03,S3aH1,S1b0205h05,S1b0205h01,S1b0205h20,S2c0415h03,S2c0415h05,S2c0415h24,S1b0205h03,S2aH8,S2d0427h20,S1b0205h24,S2c0415h22,S2d0427h05,S2d0427h22,S3aH8"

$wb.Application.Calculate()
